$wb = $excel.ActiveWorkbook

# --- Summary sheet updates ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.74   # Current Capital
$summary.Range("B4").Value = -0.26     # Total P&L $
$summary.Range("B5").Value = -0.14     # Total P&L %
$summary.Range("B6").Value = 38        # Total Trades
$summary.Range("B8").Value = 19        # Losing Trades
$summary.Range("B9").Value = 28.95     # Win Rate %

# --- Strategy Status sheet updates (MarketMaking row) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.73999999999999
$status.Range("D4").Value = 38
$status.Range("E4").Value = -0.26
$status.Range("F4").Value = -0.26
$status.Range("G4").Value = 28.95

# --- Add new trade row (#38) to both "All Trades" and "MarketMaking" sheets ---
$sheetNames = @("All Trades", "MarketMaking")
foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 39

    $ws.Cells.Item($row, 1).Value = 38

    $ws.Cells.Item($row, 2).Value = "'2026-02-17"
    $ws.Cells.Item($row, 2).ClearFormats()

    $ws.Cells.Item($row, 3).Value = "'15:23:29"
    $ws.Cells.Item($row, 3).ClearFormats()

    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "DOWN"
    $ws.Cells.Item($row, 6).Value = 0.6
    $ws.Cells.Item($row, 7).Value = 0.58
    $ws.Cells.Item($row, 8).Value = "CLOSED"
    $ws.Cells.Item($row, 9).Value = -3.3333
    $ws.Cells.Item($row, 10).Value = -0.02
    $ws.Cells.Item($row, 11).Value = 99.73999999999999
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.14
}
